# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value  = 388
$wsExhibit.Range("F9").Value  = 559
$wsExhibit.Range("F10").Value = 9
$wsExhibit.Range("F13").Value = 13570
$wsExhibit.Range("F17").Value = 5575
$wsExhibit.Range("F19").Value = 64

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F22").Value = 388
$wsAll.Range("F31").Value = 559
$wsAll.Range("F32").Value = 9
$wsAll.Range("F35").Value = 13571
$wsAll.Range("F40").Value = 5575
$wsAll.Range("F42").Value = 64
